# Apply the "cryptos list" price/volume refresh described by the commit diff.
# For D-column values that parse as plain numbers, a leading apostrophe forces
# Excel to store them as text (matching the source data, which are formatted
# strings like "556.66", not numeric cells) instead of auto-converting to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '60.030.81'
$ws.Range("E2").Value = '  +0.38%  '
# Row 3: Ethereum
$ws.Range("D3").Value = '2.410.13'
$ws.Range("E3").Value = '  -0.51%  '
# Row 4: TetherUSD
$ws.Range("E4").Value = '  -0.03%  '
# Row 5: BNB
$ws.Range("D5").Value = '''556.66'
$ws.Range("E5").Value = '  +0.80%  '
# Row 6: Solana
$ws.Range("D6").Value = '''135.74'
$ws.Range("E6").Value = '  -1.05%  '
# Row 7: USDC
$ws.Range("E7").Value = '  -0.01%  '
# Row 8: XRP
$ws.Range("D8").Value = '''0.586'
$ws.Range("E8").Value = '  -0.37%  '
# Row 9: Dogecoin
$ws.Range("E9").Value = '  -0.51%  '
# Row 10: Toncoin
$ws.Range("D10").Value = '''5.62'
$ws.Range("E10").Value = '  -1.16%  '
# Row 11: TRON
$ws.Range("E11").Value = '  -0.09%  '
# Row 12: Cardano
$ws.Range("D12").Value = '''0.349'
$ws.Range("E12").Value = '  -1.27%  '
# Row 13: Avalanche
$ws.Range("D13").Value = '''24.71'
$ws.Range("E13").Value = '  -0.33%  '
# Row 14: WrappedliquidstakedEther2.0
$ws.Range("D14").Value = '2.838.83'
$ws.Range("E14").Value = '  -0.50%  '
# Row 15: WrappedBTC
$ws.Range("D15").Value = '59.934.42'
$ws.Range("E15").Value = '  +0.32%  '
# Row 16: ShibaInu
$ws.Range("D16").Value = '''0.0000138'
$ws.Range("E16").Value = '  +0.38%  '
# Row 17: WrappedEther
$ws.Range("D17").Value = '2.414.45'
$ws.Range("E17").Value = '  -0.13%  '
# Row 18: Chainlink
$ws.Range("D18").Value = '''11.18'
$ws.Range("E18").Value = '  -0.80%  '
# Row 19: Polkadot
$ws.Range("D19").Value = '''4.51'
$ws.Range("E19").Value = '  +2.91%  '
# Row 20: BitcoinCash
$ws.Range("D20").Value = '''327.70'
$ws.Range("E20").Value = '  -0.90%  '
# Row 21: Uniswap
$ws.Range("D21").Value = '''6.78'
$ws.Range("E21").Value = '  +1.49%  '
# Row 22: Dai
$ws.Range("E22").Value = '  -0.25%  '
# Row 23: Litecoin
$ws.Range("D23").Value = '''64.69'
$ws.Range("E23").Value = '  -1.55%  '
# Row 24: Kaspa
$ws.Range("D24").Value = '''0.176'
$ws.Range("E24").Value = '  +2.70%  '
# Row 25: InternetComputer(DFINITY)
$ws.Range("D25").Value = '''8.58'
$ws.Range("E25").Value = '  -0.68%  '
# Row 26: Binance-PegBSC-USD
$ws.Range("E26").Value = '  -0.17%  '
# Row 27: Fetch.AI
$ws.Range("E27").Value = '  +3.23%  '
# Row 28: PancakeSwap
$ws.Range("D28").Value = '''1.80'
$ws.Range("E28").Value = '  +1.74%  '
# Row 29: PEPE
$ws.Range("D29").Value = '0.0₃0769'
$ws.Range("E29").Value = '  -0.92%  '
# Row 30: Monero
$ws.Range("D30").Value = '''169.94'
# Row 31: Aptos
$ws.Range("E31").Value = '  -0.36%  '
# Row 32: SuiNetwork
$ws.Range("E32").Value = '  +7.38%  '
# Row 33: PolygonEcosystemToken
$ws.Range("D33").Value = '''0.401'
$ws.Range("E33").Value = '  -2.02%  '
# Row 34: EthereumClassic
$ws.Range("D34").Value = '''18.42'
$ws.Range("E34").Value = '  -1.36%  '
# Row 35: USDe
$ws.Range("E35").Value = '  +0.09%  '
# Row 36: ImmutableX
$ws.Range("D36").Value = '''1.33'
$ws.Range("E36").Value = '  +2.93%  '
# Row 37: FirstDigitalUSD
$ws.Range("E37").Value = '  +0.06%  '
# Row 38: NEARProtocol
$ws.Range("D38").Value = '''4.20'
$ws.Range("E38").Value = '  +0.07%  '
# Row 39: Bittensor
$ws.Range("D39").Value = '''323.19'
$ws.Range("E39").Value = '  +3.14%  '
# Row 40: Stacks
$ws.Range("D40").Value = '''1.60'
$ws.Range("E40").Value = '  -0.53%  '
# Row 41: OKB
$ws.Range("D41").Value = '''38.45'
$ws.Range("E41").Value = '  -2.41%  '
# Row 42: Aave
$ws.Range("D42").Value = '''147.07'
$ws.Range("E42").Value = '  +6.13%  '
# Row 43: Filecoin
$ws.Range("D43").Value = '''3.59'
$ws.Range("E43").Value = '  -1.98%  '
# Row 44: Stellar
$ws.Range("D44").Value = '''0.0969'
$ws.Range("E44").Value = '  +0.09%  '
# Row 45: InjectiveProtocol
$ws.Range("D45").Value = '''19.93'
$ws.Range("E45").Value = '  +2.13%  '
# Row 46: Hedera
$ws.Range("D46").Value = '''0.0516'
$ws.Range("E46").Value = '  -0.22%  '
# Row 47: Mantle
$ws.Range("D47").Value = '''0.576'
$ws.Range("E47").Value = '  -0.83%  '
# Row 48: VeChain
$ws.Range("E48").Value = '  -1.68%  '
# Row 49: WhiteBITCoin
$ws.Range("E49").Value = '  -0.18%  '
# Row 50: dogwifhat
$ws.Range("D50").Value = '''1.58'
$ws.Range("E50").Value = '  -0.46%  '
# Row 51: ZEEBU
$ws.Range("E51").Value = '  -0.66%  '
